# 1. 加入MyUser Model
#    Add a new worksheet "MyUser" at the end of the workbook, mirroring the
#    layout used by the other "table schema" sheets (e.g. "User"), and make
#    it the active/selected sheet (replacing " Problem" as the tab that was
#    previously selected).

$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "MyUser"

$newSheet.Range("A1").Value = "MyUser"
$newSheet.Range("B1").Value = "自加字段 OneToOne - User"

$newSheet.Range("A2").Value = "id"
$newSheet.Range("B2").Value = "auto"

$newSheet.Range("A3").Value = "nickname"
$newSheet.Range("B3").Value = "nvarchar"
$newSheet.Range("C3").Value = "姓名"

$newSheet.Range("A4").Value = "authority"
$newSheet.Range("B4").Value = "int"
$newSheet.Range("C4").Value = "权限 0学生1教师2管理员"

$newSheet.Range("A5").Value = "user_id"
$newSheet.Range("B5").Value = "int"

# Match the column widths ("best fit" around the longest cell text) used on
# the sibling schema sheets.
$newSheet.Columns.Item(1).ColumnWidth = 9.714285714285714
$newSheet.Columns.Item(3).ColumnWidth = 23

# Put the cursor on C5 (the last populated cell), matching the saved
# selection, and make this newly added sheet the active tab - it becomes the
# one shown/selected when the workbook is reopened.
$newSheet.Range("C5").Select() | Out-Null
$newSheet.Activate()
